$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.655.13'
$ws.Range('E2').Value = '  -6.95%  '
$ws.Range('D3').Value = '1.699.44'
$ws.Range('E3').Value = '  -5.52%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '219.52'
$ws.Range('E5').Value = '  -5.18%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.5120'
$ws.Range('E6').Value = '  -13.36%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.004'
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2651'
$ws.Range('E8').Value = '  -4.26%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '22.18'
$ws.Range('E9').Value = '  -4.49%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06305'
$ws.Range('E10').Value = '  -7.37%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07355'
$ws.Range('E11').Value = '  -2.14%  '
$ws.Range('D12').Value = '1.697.23'
$ws.Range('E12').Value = '  -5.73%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '4.516'
$ws.Range('E13').Value = '  -5.15%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '0.5844'
$ws.Range('E14').Value = '  -5.91%  '
$ws.Range('D15').Value = '1.930.13'
$ws.Range('E15').Value = '  -5.54%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.000008428'
$ws.Range('E16').Value = '  -7.64%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '65.66'
$ws.Range('E17').Value = '  -13.06%  '
$ws.Range('D18').Value = '26.668.10'
$ws.Range('E18').Value = '  -6.83%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '5.020'
$ws.Range('E19').Value = '  -8.07%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '1.004'
$ws.Range('E20').Value = '  +0.11%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.02'
$ws.Range('E21').Value = '  -4.17%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '187.09'
$ws.Range('E22').Value = '  -11.14%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '6.274'
$ws.Range('E23').Value = '  -8.03%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.004'
$ws.Range('E24').Value = '  +0.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '144.57'
$ws.Range('E25').Value = '  -5.95%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.528'
$ws.Range('E26').Value = '  -4.11%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.1161'
$ws.Range('E27').Value = '  -8.50%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '15.72'
$ws.Range('E28').Value = '  -4.33%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.349'
$ws.Range('E29').Value = '  -4.80%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.05683'
$ws.Range('E30').Value = '  -8.19%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.339'
$ws.Range('E31').Value = '  -5.87%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.514'
$ws.Range('E32').Value = '  -6.96%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '3.491'
$ws.Range('E33').Value = '  -8.29%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.647'
$ws.Range('E34').Value = '  -4.80%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.022'
$ws.Range('E35').Value = '  -3.22%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.6046'
$ws.Range('E36').Value = '  -5.76%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.361'
$ws.Range('E37').Value = '  -5.53%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.685'
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.01614'
$ws.Range('E39').Value = '  -4.92%  '
$ws.Range('D40').Value = '1.098.36'
$ws.Range('E40').Value = '  -4.16%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.8597'
$ws.Range('E41').Value = '  -2.65%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.847'
$ws.Range('E42').Value = '  -10.69%  '
$ws.Range('E43').Value = '  -0.57%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '100.01'
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '1.856.84'
$ws.Range('E45').Value = '  -4.88%  '
$ws.Range('E46').Value = '  -2.62%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '56.82'
$ws.Range('E47').Value = '  -5.76%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '8.204'
$ws.Range('E48').Value = '  -1.75%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.004'
$ws.Range('E49').Value = '  +0.23%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.05238'
$ws.Range('E50').Value = '  -4.19%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.4319'
$ws.Range('E51').Value = '  -3.60%  '
